$p = $ppt.ActivePresentation

# --- 1) Swap the presentation's theme colours (theme1.xml) from the
#        "Integral" / Red Violet palette to the default "Office Theme"
#        palette, via the per-slide ThemeColorScheme object (this is the
#        only path that edits the live theme part without disturbing
#        other OOXML parts).
$tcs = $p.Slides.Item(1).ThemeColorScheme
$tcs.Colors(1).RGB  = 0          # dk1      000000
$tcs.Colors(2).RGB  = 16777215   # lt1      FFFFFF
$tcs.Colors(3).RGB  = 6968388    # dk2      44546A
$tcs.Colors(4).RGB  = 15132391   # lt2      E7E6E6
$tcs.Colors(5).RGB  = 13998939   # accent1  5B9BD5
$tcs.Colors(6).RGB  = 3243501    # accent2  ED7D31
$tcs.Colors(7).RGB  = 10855845   # accent3  A5A5A5
$tcs.Colors(8).RGB  = 49407      # accent4  FFC000
$tcs.Colors(9).RGB  = 12874308   # accent5  4472C4
$tcs.Colors(10).RGB = 4697456    # accent6  70AD47
$tcs.Colors(11).RGB = 12673797   # hlink    0563C1
$tcs.Colors(12).RGB = 7491477    # folHlink 954F72

# --- 2) Re-style the three summary tables (slides 14-16) with the
#        built-in "Medium Style 2 - Accent 1" table style.
for ($si = 14; $si -le 16; $si++) {
    $slide = $p.Slides.Item($si)
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $shp = $slide.Shapes.Item($i)
        if ($shp.HasTable) {
            $shp.Table.ApplyStyle("{663E054E-A488-44D7-9480-B80184B4BE35}")
        }
    }
}
